$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "target" trial rows (currently rows 13,14,15 data + row16 H/I/J/K)
# down by one row (13->14, 14->15, 15->16), then turn row 12 into the new
# catch trial row (clearing its old target-trial content), and clear the
# catch-specific cells that used to live on row 16 (M16:V16 stay filled
# since row16 becomes a target row; H16/I16 get filled; J16/K16 get the
# "old"/"j" values; L16 gets the old row15 image).

# --- Row 16 becomes what used to be row 15's target data, plus H/I/J/K filled in ---
$ws.Range("H16").Value = "bedrooms"
$ws.Range("I16").Value = "target"
$ws.Range("J16").Value = "old"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_1vq1v.png"
$ws.Range("M16").Value = 69.42857142857143
$ws.Range("N16").Value = 46.59523809523809
$ws.Range("O16").Value = 58.01190476190476
$ws.Range("P16").Value = 42
$ws.Range("Q16").Value = 5
$ws.Range("R16").Value = 5
$ws.Range("S16").Value = 5
$ws.Range("T16").Value = 5
$ws.Range("U16").Value = 5
$ws.Range("V16").Value = 5

# --- Row 15 gets what used to be row 14's target data ---
$ws.Range("L15").Value = "stimuli/img_oou46.png"
$ws.Range("M15").Value = 75.70270270270271
$ws.Range("N15").Value = 54.86486486486486
$ws.Range("O15").Value = 65.28378378378379
$ws.Range("P15").Value = 37
$ws.Range("Q15").Value = 6
$ws.Range("R15").Value = 6
$ws.Range("S15").Value = 6
$ws.Range("T15").Value = 6
$ws.Range("U15").Value = 6
$ws.Range("V15").Value = 6

# --- Row 14 gets what used to be row 13's target data ---
$ws.Range("L14").Value = "stimuli/img_a9acb.png"
$ws.Range("M14").Value = 77.11428571428571
$ws.Range("N14").Value = 58.42857142857143
$ws.Range("O14").Value = 67.77142857142857
$ws.Range("P14").Value = 35
$ws.Range("Q14").Value = 7
$ws.Range("R14").Value = 7
$ws.Range("S14").Value = 7
$ws.Range("T14").Value = 7
$ws.Range("U14").Value = 7
$ws.Range("V14").Value = 7

# --- Row 13 gets what used to be row 12's target data ---
$ws.Range("L13").Value = "stimuli/img_cogrz.png"
$ws.Range("M13").Value = 60.5
$ws.Range("N13").Value = 39.71428571428572
$ws.Range("O13").Value = 50.10714285714286
$ws.Range("P13").Value = 42
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = 3
$ws.Range("S13").Value = 3
$ws.Range("T13").Value = 3
$ws.Range("U13").Value = 3
$ws.Range("V13").Value = 3

# --- Row 12 becomes the new catch trial row; clear old target-only cells ---
$ws.Range("H12").ClearContents()
$ws.Range("I12").ClearContents()
$ws.Range("J12").Value = "catch"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/catch_20.jpg"
$ws.Range("M12:V12").ClearContents()
